# Add a new weekly price-report group (3 rows) for "Vega Monumental Concepción - Sandia"
# at the top of the data (row 27), pushing all existing records down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 27 (existing rows 27-118 become 30-121)
$ws.Rows("27:29").Insert()

# Row 27: new "Extra" quality record for 2022-12-06 (serial 44901)
$ws.Range("A27").Value2 = 11
$ws.Range("B27").Value2 = "Vega Monumental Concepción"
$ws.Range("C27").Value2 = "Bíobío"
$ws.Range("D27").Value2 = 44901
$ws.Range("E27").Value2 = 8
$ws.Range("F27").Value2 = 100112028
$ws.Range("G27").Value2 = "Sandia"
$ws.Range("H27").Value2 = "Sin especificar"
$ws.Range("I27").Value2 = "Extra"
$ws.Range("J27").Value2 = 2700
$ws.Range("K27").Value2 = 3800
$ws.Range("L27").Value2 = 4000
$ws.Range("M27").Value2 = 3889
$ws.Range("N27").Value2 = "`$/unidad"
$ws.Range("O27").Value2 = "Región de O'Higgins"
$ws.Range("P27").Value2 = 3889
$ws.Range("Q27").Value2 = 1
$ws.Range("R27").Value2 = "Hortaliza"

# Row 28: new "Primera" quality record for 2022-12-06 (serial 44901)
$ws.Range("A28").Value2 = 11
$ws.Range("B28").Value2 = "Vega Monumental Concepción"
$ws.Range("C28").Value2 = "Bíobío"
$ws.Range("D28").Value2 = 44901
$ws.Range("E28").Value2 = 8
$ws.Range("F28").Value2 = 100112028
$ws.Range("G28").Value2 = "Sandia"
$ws.Range("H28").Value2 = "Sin especificar"
$ws.Range("I28").Value2 = "Primera"
$ws.Range("J28").Value2 = 2200
$ws.Range("K28").Value2 = 3200
$ws.Range("L28").Value2 = 3500
$ws.Range("M28").Value2 = 3364
$ws.Range("N28").Value2 = "`$/unidad"
$ws.Range("O28").Value2 = "Región de O'Higgins"
$ws.Range("P28").Value2 = 3364
$ws.Range("Q28").Value2 = 1
$ws.Range("R28").Value2 = "Hortaliza"

# Row 29: new "Segunda" quality record for 2022-12-06 (serial 44901)
$ws.Range("A29").Value2 = 11
$ws.Range("B29").Value2 = "Vega Monumental Concepción"
$ws.Range("C29").Value2 = "Bíobío"
$ws.Range("D29").Value2 = 44901
$ws.Range("E29").Value2 = 8
$ws.Range("F29").Value2 = 100112028
$ws.Range("G29").Value2 = "Sandia"
$ws.Range("H29").Value2 = "Sin especificar"
$ws.Range("I29").Value2 = "Segunda"
$ws.Range("J29").Value2 = 2000
$ws.Range("K29").Value2 = 2700
$ws.Range("L29").Value2 = 3000
$ws.Range("M29").Value2 = 2850
$ws.Range("N29").Value2 = "`$/unidad"
$ws.Range("O29").Value2 = "Región de O'Higgins"
$ws.Range("P29").Value2 = 2850
$ws.Range("Q29").Value2 = 1
$ws.Range("R29").Value2 = "Hortaliza"
